$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item(1)

# Add a new worksheet right after Sheet1 -> becomes "Sheet2"
$ws2 = $wb.Worksheets.Add($null, $ws1)

# Header row
$ws2.Range("C3").Value = "Name"
$ws2.Range("D3").Value = "Salary "
$ws2.Range("E3").Value = "Workexp"

# Data rows
$ws2.Range("C4").Value = "Ram"
$ws2.Range("D4").Value = 30000
$ws2.Range("E4").Value = 5

$ws2.Range("C5").Value = "Shyam"
$ws2.Range("D5").Value = 25000
$ws2.Range("E5").Value = 6

$ws2.Range("C6").Value = "Geeta"
$ws2.Range("D6").Value = 35000
$ws2.Range("E6").Value = 4

$ws2.Range("C7").Value = "Malini"
$ws2.Range("D7").Value = 32000
$ws2.Range("E7").Value = 3

# Labels for the lookup demo
$ws2.Range("C10").Value = "index"
$ws2.Range("E10").Value = "match"

# Formulas
$ws2.Range("C11").Formula = "=INDEX(C3:E7,4,3)"
$ws2.Range("E11").Formula = "=MATCH(C4,C3:C7,0)"

# Match the selection/view state recorded in the target file
$ws2.Range("F16").Select()
